$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.679.60"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.944.12"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.543"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0845"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").Value = "3.404.83"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "2.934.75"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.955"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.95%  "
$ws.Range("D18").Value = "51.637.22"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +20.32%  "
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("E28").Value = "  -3.52%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "52.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "34.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  -6.04%  "
$ws.Range("E42").Value = "  -4.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.279"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.76%  "
$ws.Range("D48").Value = "2.031.86"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0326"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.59%  "
